# Update content of toolbox and descriptions
# Rewrites columns A-H for rows 2-22 of the orthomyxoviruses "family" sheet
# to reflect the corrected/updated virus toolbox data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('H3N2', 'H3N2', 'A/H3N2', 'Flu', 'Influenza (Flu)', 'Influenza A', 'Yes', 'Yes'),
    @('H7N9', 'H7N9', 'A/H7N9', 'AI', 'Avian influenza (Bird Flu)', 'Influenza A', 'Yes', 'Yes'),
    @('Victoria', 'B/Victoria', 'B/Victoria', 'Flu', 'Influenza (Flu)', 'Influenza B', 'Yes', 'Yes'),
    @('Yamagata', 'B/Yamagata', 'B/Yamagata', 'Flu', 'Influenza (Flu)', 'Influenza B', 'Yes', 'Yes'),
    @('H1N1', 'H1N1', 'A/H1N1', 'Flu', 'Influenza (Flu)', 'Influenza A', 'Yes', 'Yes'),
    @('H5N1', 'H5N1', 'A/H5N1', 'AI', 'Avian influenza (Bird Flu)', 'Influenza A', 'Yes', 'Yes'),
    @('H9N2', 'H9N2', 'A/H9N2', 'AI', 'Avian influenza (Bird Flu)', 'Influenza A', 'Yes', 'Yes'),
    @('H3N8', 'H3N8', 'A/H3N8', 'AI', 'Avian influenza (Bird Flu)', 'Influenza A', 'Yes', 'Yes'),
    @('H4N6', 'H4N6', 'A/H4N6', '', '', 'Influenza A', 'Yes', 'No'),
    @('H14N3', 'H14N3', 'A/H14N3', '', '', 'Influenza A', 'Yes', 'No'),
    @('H15N9', 'H15N9', 'A/H15N9', '', '', 'Influenza A', 'Yes', 'No'),
    @('H10N7', 'H10N7', 'A/H10N7', '', '', 'Influenza A', 'Yes', 'No'),
    @('H2N2', 'H2N2', 'A/H2N2', '', '', 'Influenza A', 'Yes', 'No'),
    @('H6N2', 'H6N2', 'A/H6N2', '', '', 'Influenza A', 'Yes', 'No'),
    @('H17N10', 'H17N10', 'A/H17N10', '', '', 'Influenza A', 'Yes', 'No'),
    @('H18N11', 'H18N11', 'A/H18N11', '', '', 'Influenza A', 'Yes', 'No'),
    @('H8N4', 'H8N4', 'A/H8N4', '', '', 'Influenza A', 'Yes', 'No'),
    @('H11N9', 'H11N9', 'A/H11N9', '', '', 'Influenza A', 'Yes', 'No'),
    @('H12N5', 'H12N5', 'A/H12N5', '', '', 'Influenza A', 'Yes', 'No'),
    @('H13N6', 'H13N6', 'A/H13N6', '', '', 'Influenza A', 'Yes', 'No'),
    @('H16N3', 'H16N3', 'A/H16N3', '', '', 'Influenza A', 'Yes', 'No')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $colLetter = [char](65 + $j)
        $ws.Range("$colLetter$rowNum").Value = $rowValues[$j]
    }
}
